$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rentometer")

# Update quickview_url (B17) and its hyperlink with the new token
$ws.Range("B17").Value = "https://www.rentometer.com/analysis/3-bed/317-newell-st-barberton-oh-44203/5huWv3zgogk/quickview"

# Update credits_remaining (B18)
$ws.Range("B18").Value = 1943

# Update token (B19)
$ws.Range("B19").Value = "5huWv3zgogk"

# Update links (B20)
$ws.Range("B20").Value = "[{'rel': 'request pro report', 'href': 'https://www.rentometer.com/api/v1/request_pro_report?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=5huWv3zgogk'}, {'rel': 'nearby comps', 'href': 'https://www.rentometer.com/api/v1/nearby_comps?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=5huWv3zgogk'}]"
